$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = "변수분리법"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/06/separable_differential_equations.html"

# Row 8
$ws.Range("D8").Value = "카카오브레인"

# Row 20
$ws.Range("D20").Value = "[AI] 간단한 딥러닝 웹서비스 - 마스크 착용 감지 프로젝트 (w/ Teachable Machine)"

# Row 26
$ws.Range("D26").Value = "bayesian optimization in trading"
$ws.Range("E26").Value = "https://blog.est.ai/2021/05/bayesian-optimization-in-trading/"

# Row 32
$ws.Range("D32").Value = "Kernel Density Estimation (KDE)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/318"
